$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D price cells to stay as text (many values look like plain
# decimal numbers, e.g. "594.26", and Excel would otherwise silently coerce
# them to numeric values on assignment). Temporarily mark the whole Price
# column as Text, write the new values, then clear the format again so the
# cells end up with no explicit style (matching the original inlineStr
# cells, which carry no style index).
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

# Column D (Price) updates
$ws.Range("D2").Value = "73.379.86"
$ws.Range("D3").Value = "4.056.75"
$ws.Range("D5").Value = "594.26"
$ws.Range("D6").Value = "153.80"
$ws.Range("D7").Value = "4.051.45"
$ws.Range("D12").Value = "53.60"
$ws.Range("D14").Value = "11.07"
$ws.Range("D15").Value = "4.704.09"
$ws.Range("D16").Value = "4.053.94"
$ws.Range("D17").Value = "14.29"
$ws.Range("D18").Value = "1.24"
$ws.Range("D19").Value = "20.79"
$ws.Range("D20").Value = "73.281.44"
$ws.Range("D22").Value = "444.27"
$ws.Range("D23").Value = "4.71"
$ws.Range("D24").Value = "97.62"
$ws.Range("D27").Value = "4.33"
$ws.Range("D28").Value = "11.44"
$ws.Range("D31").Value = "37.01"
$ws.Range("D32").Value = "7.91"
$ws.Range("D34").Value = "13.67"
$ws.Range("D35").Value = "690.78"
$ws.Range("D36").Value = "48.49"
$ws.Range("D37").Value = "68.85"
$ws.Range("D38").Value = "0.448"
$ws.Range("D39").Value = "0.0₃0879"
$ws.Range("D46").Value = "1.00"
$ws.Range("D49").Value = "3.39"
$ws.Range("D50").Value = "3.54"
$ws.Range("D43").Value = "3.36"

# Drop the temporary Text format so the cells end up unstyled again.
$priceRange.ClearFormats()

# Column E (Volume 1h) updates - plain percentage-like text, safe to assign directly
$ws.Range("E2").Value = "  +1.60%  "
$ws.Range("E3").Value = "  +0.66%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("E5").Value = "  +12.22%  "
$ws.Range("E6").Value = "  +1.66%  "
$ws.Range("E7").Value = "  +0.44%  "
$ws.Range("E8").Value = "  -1.02%  "
$ws.Range("E9").Value = "  +0.01%  "
$ws.Range("E10").Value = "  +1.53%  "
$ws.Range("E11").Value = "  -0.73%  "
$ws.Range("E12").Value = "  +11.90%  "
$ws.Range("E13").Value = "  -0.94%  "
$ws.Range("E14").Value = "  +3.62%  "
$ws.Range("E15").Value = "  +0.87%  "
$ws.Range("E16").Value = "  +0.95%  "
$ws.Range("E17").Value = "  +0.77%  "
$ws.Range("E18").Value = "  +3.64%  "
$ws.Range("E19").Value = "  +0.66%  "
$ws.Range("E20").Value = "  +1.66%  "
$ws.Range("E21").Value = "  -0.68%  "
$ws.Range("E22").Value = "  +3.78%  "
$ws.Range("E23").Value = "  +12.01%  "
$ws.Range("E24").Value = "  -0.64%  "
$ws.Range("E25").Value = "  +1.07%  "
$ws.Range("E26").Value = "  +1.23%  "
$ws.Range("E27").Value = "  +20.55%  "
$ws.Range("E28").Value = "  +1.78%  "
$ws.Range("E29").Value = "  +1.40%  "
$ws.Range("E30").Value = "  +2.07%  "
$ws.Range("E31").Value = "  +0.63%  "
$ws.Range("E32").Value = "  +10.20%  "
$ws.Range("E33").Value = "  +3.87%  "
$ws.Range("E34").Value = "  +1.90%  "
$ws.Range("E35").Value = "  +2.07%  "
$ws.Range("E36").Value = "  +8.61%  "
$ws.Range("E37").Value = "  +4.10%  "
$ws.Range("E38").Value = "  -0.54%  "
$ws.Range("E39").Value = "  +5.81%  "
$ws.Range("E40").Value = "  +16.80%  "
$ws.Range("E41").Value = "  -2.65%  "
$ws.Range("E42").Value = "  +5.09%  "
$ws.Range("E43").Value = "  -0.84%  "
$ws.Range("E44").Value = "  -0.05%  "
$ws.Range("E45").Value = "  +2.01%  "
$ws.Range("E46").Value = "  +0.11%  "
$ws.Range("E47").Value = "  +0.64%  "
$ws.Range("E48").Value = "  +2.45%  "
$ws.Range("E49").Value = "  -2.31%  "
$ws.Range("E50").Value = "  +7.37%  "
$ws.Range("E51").Value = "  +2.13%  "

# Rows 42 and 43 swap positions for Coin (B) and Link (C): ThetaToken <-> WEMIXToken
$ws.Range("B42").Value = "WEMIXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("B43").Value = "ThetaToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
